$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = ""
$ws.Range("H87").Value = 33354
$ws.Range("J87").Value = 33354
$ws.Range("L87").Value = 33354
$ws.Range("N87").Value = -35850
$ws.Range("H90").Value = 33354
$ws.Range("J90").Value = 33354
$ws.Range("L90").Value = 100062
$ws.Range("N90").Value = -112542
$ws.Range("H112").Value = 4615.3335
$ws.Range("J112").Value = 4694.615
$ws.Range("L112").Value = 14083.845
$ws.Range("N112").Value = -16299.845
$ws.Range("H132").Value = 2556.8
$ws.Range("I132").Value = 1245.0952
$ws.Range("K132").Value = 3735.2856
$ws.Range("M132").Value = -1205.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28316.96
$ws.Range("I32").Value = 15583.057
$ws.Range("K32").Value = 15583.057
$ws.Range("M32").Value = -15296.057
$ws.Range("H61").Value = 1158.2565
$ws.Range("I61").Value = 780.1
$ws.Range("K61").Value = 780.1
$ws.Range("M61").Value = -568.1
$ws.Range("H136").Value = 1158.2565
$ws.Range("I136").Value = 780.1
$ws.Range("K136").Value = 2340.3
$ws.Range("M136").Value = 209.6999999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 278
$ws.Range("I22").Value = 278
$ws.Range("K22").Value = 278
$ws.Range("M22").Value = -105
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = ""
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 1999
$ws.Range("I23").Value = 1999
$ws.Range("K23").Value = 1999
$ws.Range("M23").Value = -1759
$ws.Range("H27").Value = 1999
$ws.Range("I27").Value = 1999
$ws.Range("K27").Value = 1999
$ws.Range("M27").Value = -1807
$ws.Range("H31").Value = 2162.1538
$ws.Range("I31").Value = 1522.52
$ws.Range("J31").Value = 3304.3572
$ws.Range("K31").Value = 1522.52
$ws.Range("L31").Value = 3304.3572
$ws.Range("M31").Value = -1227.52
$ws.Range("N31").Value = -3894.3572
$ws.Range("H34").Value = 2162.1538
$ws.Range("I34").Value = 1522.52
$ws.Range("J34").Value = 3304.3572
$ws.Range("K34").Value = 1522.52
$ws.Range("L34").Value = 3304.3572
$ws.Range("M34").Value = -1320.52
$ws.Range("N34").Value = -3708.3572
$ws.Range("H47").Value = 37000
$ws.Range("J47").Value = 37000
$ws.Range("L47").Value = 37000
$ws.Range("N47").Value = -38132
$ws.Range("H99").Value = 71998.5
$ws.Range("I99").Value = 129999
$ws.Range("K99").Value = 129999
$ws.Range("M99").Value = -128501
$ws.Range("H126").Value = 71998.5
$ws.Range("I126").Value = 129999
$ws.Range("K126").Value = 389997
$ws.Range("M126").Value = -387527
$ws.Range("H134").Value = 2214.9736
$ws.Range("I134").Value = 2242.8125
$ws.Range("K134").Value = 6728.4375
$ws.Range("M134").Value = -4193.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1332.0625
$ws.Range("I34").Value = 165.28572
$ws.Range("J34").Value = 2239.5557
$ws.Range("K34").Value = 495.85716
$ws.Range("L34").Value = 6718.6671
$ws.Range("M34").Value = -411.85716
$ws.Range("N34").Value = -6886.6671
$ws.Range("H39").Value = 3305.2222
$ws.Range("J39").Value = 3750
$ws.Range("L39").Value = 11250
$ws.Range("N39").Value = -11838
$ws.Range("H41").Value = 276
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = ""
$ws.Range("H48").Value = 1999.6666
$ws.Range("I48").Value = 1999
$ws.Range("K48").Value = 5997
$ws.Range("M48").Value = -5747
$ws.Range("H55").Value = 10419379
$ws.Range("J55").Value = 12503050
$ws.Range("L55").Value = 37509150
$ws.Range("N55").Value = -37509504
$ws.Range("H64").Value = 4112.375
$ws.Range("I64").Value = 1449.5
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 4348.5
$ws.Range("L64").Value = 15000
$ws.Range("M64").Value = -4078.5
$ws.Range("N64").Value = -15540
$ws.Range("H67").Value = 4112.375
$ws.Range("I67").Value = 1449.5
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 4348.5
$ws.Range("L67").Value = 15000
$ws.Range("M67").Value = -3412.5
$ws.Range("N67").Value = -16872
$ws.Range("H105").Value = 10833.333
$ws.Range("J105").Value = 15000
$ws.Range("L105").Value = 45000
$ws.Range("N105").Value = -50242
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = ""
$ws.Range("H128").Value = 399686.75
$ws.Range("I128").Value = 399686.75
$ws.Range("K128").Value = 1199060.25
$ws.Range("M128").Value = -1194080.25
$ws.Range("H134").Value = 6943.1665
$ws.Range("I134").Value = 2914.75
$ws.Range("K134").Value = 8744.25
$ws.Range("M134").Value = -3674.25
$ws.Range("H136").Value = 6999.4443
$ws.Range("I136").Value = 2748.75
$ws.Range("K136").Value = 8246.25
$ws.Range("M136").Value = -3146.25
$ws.Range("H138").Value = 1502.5
$ws.Range("I138").Value = 996.4
$ws.Range("J138").Value = 4033
$ws.Range("K138").Value = 2989.2
$ws.Range("L138").Value = 12099
$ws.Range("M138").Value = 2150.8
$ws.Range("N138").Value = -22379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 30000
$ws.Range("J40").Value = 30000
$ws.Range("L40").Value = 30000
$ws.Range("N40").Value = -30302
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H92").Value = 11250
$ws.Range("J92").Value = 11250
$ws.Range("L92").Value = 11250
$ws.Range("N92").Value = -14994
$ws.Range("H136").Value = 45888.375
$ws.Range("J136").Value = 45888.375
$ws.Range("L136").Value = 137665.125
$ws.Range("N136").Value = -142765.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 122.5
$ws.Range("J2").Value = 122.5
$ws.Range("L2").Value = 122.5
$ws.Range("N2").Value = -346.5
$ws.Range("H22").Value = 1248
$ws.Range("J22").Value = 1686.5
$ws.Range("L22").Value = 1686.5
$ws.Range("N22").Value = -2276.5
$ws.Range("H27").Value = 1248
$ws.Range("J27").Value = 1686.5
$ws.Range("L27").Value = 1686.5
$ws.Range("N27").Value = -1900.5
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = ""
$ws.Range("H81").Value = 4094.6667
$ws.Range("J81").Value = 740
$ws.Range("L81").Value = 1480
$ws.Range("N81").Value = -3602
$ws.Range("H84").Value = 4094.6667
$ws.Range("J84").Value = 740
$ws.Range("L84").Value = 7400
$ws.Range("N84").Value = -18008
$ws.Range("H136").Value = 4343.2856
$ws.Range("I136").Value = 4777.1177
$ws.Range("J136").Value = 2499.5
$ws.Range("K136").Value = 14331.3531
$ws.Range("L136").Value = 7498.5
$ws.Range("M136").Value = -11781.3531
$ws.Range("N136").Value = -12598.5
